$d = $word.ActiveDocument

# Update the date line at the top of the document (first paragraph).
# Assigning directly to the paragraph Range.Text preserves run formatting
# (font/size) while swapping only the text content.
$d.Paragraphs.Item(1).Range.Text = "2023-03-25 Saturday"

# Update each math expression in the table by absolute cell position.
# NOTE: this runtime's Find/Execute(Replace:=wdReplaceAll) operates on the
# whole document regardless of which Range's Find object is invoked, so
# scoped Find-and-replace is unsafe here (several old/new values collide,
# e.g. "5+65=" -> "61+1=" while a later cell holds old value "61+1=").
# Setting Cell.Range.Text directly targets exactly one cell and keeps
# existing run formatting intact.
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "77-42="
$table.Cell(1, 2).Range.Text = "76-53="
$table.Cell(1, 3).Range.Text = "36+36="
$table.Cell(1, 4).Range.Text = "97-59="
$table.Cell(1, 5).Range.Text = "91-91="
$table.Cell(2, 1).Range.Text = "55+1="
$table.Cell(2, 2).Range.Text = "30-9="
$table.Cell(2, 3).Range.Text = "86-20="
$table.Cell(2, 4).Range.Text = "17+21="
$table.Cell(2, 5).Range.Text = "1+53="
$table.Cell(3, 1).Range.Text = "76-66="
$table.Cell(3, 2).Range.Text = "98-53="
$table.Cell(3, 3).Range.Text = "18+67="
$table.Cell(3, 4).Range.Text = "94-64="
$table.Cell(3, 5).Range.Text = "46+43="
$table.Cell(4, 1).Range.Text = "12+86="
$table.Cell(4, 2).Range.Text = "89-13="
$table.Cell(4, 3).Range.Text = "11+86="
$table.Cell(4, 4).Range.Text = "2+73="
$table.Cell(4, 5).Range.Text = "29-25="
$table.Cell(5, 1).Range.Text = "1+27="
$table.Cell(5, 2).Range.Text = "39+59="
$table.Cell(5, 3).Range.Text = "98-47="
$table.Cell(5, 4).Range.Text = "51-3="
$table.Cell(5, 5).Range.Text = "49-24="
$table.Cell(6, 1).Range.Text = "61+1="
$table.Cell(6, 2).Range.Text = "31-22="
$table.Cell(6, 3).Range.Text = "2+22="
$table.Cell(6, 4).Range.Text = "19+67="
$table.Cell(6, 5).Range.Text = "44-29="
$table.Cell(7, 1).Range.Text = "10+0="
$table.Cell(7, 2).Range.Text = "3+87="
$table.Cell(7, 3).Range.Text = "85-52="
$table.Cell(7, 4).Range.Text = "43-12="
$table.Cell(7, 5).Range.Text = "19-8="
$table.Cell(8, 1).Range.Text = "9-2="
$table.Cell(8, 2).Range.Text = "60+15="
$table.Cell(8, 3).Range.Text = "1+46="
$table.Cell(8, 4).Range.Text = "64+6="
$table.Cell(8, 5).Range.Text = "7+56="
$table.Cell(9, 1).Range.Text = "0+96="
$table.Cell(9, 2).Range.Text = "58+34="
$table.Cell(9, 3).Range.Text = "85-48="
$table.Cell(9, 4).Range.Text = "3+24="
$table.Cell(9, 5).Range.Text = "27+6="
$table.Cell(10, 1).Range.Text = "12+76="
$table.Cell(10, 2).Range.Text = "10+33="
$table.Cell(10, 3).Range.Text = "7+26="
$table.Cell(10, 4).Range.Text = "63-5="
$table.Cell(10, 5).Range.Text = "41+8="
$table.Cell(11, 1).Range.Text = "93-61="
$table.Cell(11, 2).Range.Text = "76-27="
$table.Cell(11, 3).Range.Text = "96-38="
$table.Cell(11, 4).Range.Text = "66-58="
$table.Cell(11, 5).Range.Text = "49+21="
$table.Cell(12, 1).Range.Text = "34+43="
$table.Cell(12, 2).Range.Text = "71+9="
$table.Cell(12, 3).Range.Text = "45-6="
$table.Cell(12, 4).Range.Text = "33-31="
$table.Cell(12, 5).Range.Text = "98-65="
$table.Cell(13, 1).Range.Text = "55+4="
$table.Cell(13, 2).Range.Text = "51-25="
$table.Cell(13, 3).Range.Text = "83-66="
$table.Cell(13, 4).Range.Text = "87-73="
$table.Cell(13, 5).Range.Text = "93-20="
$table.Cell(14, 1).Range.Text = "64-20="
$table.Cell(14, 2).Range.Text = "92-25="
$table.Cell(14, 3).Range.Text = "63-1="
$table.Cell(14, 4).Range.Text = "61-24="
$table.Cell(14, 5).Range.Text = "41+35="
$table.Cell(15, 1).Range.Text = "77-59="
$table.Cell(15, 2).Range.Text = "92-23="
$table.Cell(15, 3).Range.Text = "27+52="
$table.Cell(15, 4).Range.Text = "1+48="
$table.Cell(15, 5).Range.Text = "46+5="
$table.Cell(16, 1).Range.Text = "28-18="
$table.Cell(16, 2).Range.Text = "77-3="
$table.Cell(16, 3).Range.Text = "19+71="
$table.Cell(16, 4).Range.Text = "85+3="
$table.Cell(16, 5).Range.Text = "22-17="
$table.Cell(17, 1).Range.Text = "49-14="
$table.Cell(17, 2).Range.Text = "54-25="
$table.Cell(17, 3).Range.Text = "77-20="
$table.Cell(17, 4).Range.Text = "40-34="
$table.Cell(17, 5).Range.Text = "63+32="
$table.Cell(18, 1).Range.Text = "15+78="
$table.Cell(18, 2).Range.Text = "34+5="
$table.Cell(18, 3).Range.Text = "1+84="
$table.Cell(18, 4).Range.Text = "7+1="
$table.Cell(18, 5).Range.Text = "37-36="
$table.Cell(19, 1).Range.Text = "45-11="
$table.Cell(19, 2).Range.Text = "90-70="
$table.Cell(19, 3).Range.Text = "51+10="
$table.Cell(19, 4).Range.Text = "49-25="
$table.Cell(19, 5).Range.Text = "83-32="
$table.Cell(20, 1).Range.Text = "27+22="
$table.Cell(20, 2).Range.Text = "80+6="
$table.Cell(20, 3).Range.Text = "19+63="
$table.Cell(20, 4).Range.Text = "7+72="
$table.Cell(20, 5).Range.Text = "58-22="
